# Insert a new weekly record at the top of the data block (row 31),
# pushing all existing data rows (31-115) down by one (to 32-116).
# The dimension grows from A1:R115 to A1:R116.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31..115 down to 32..116, leaving a blank row 31
# (formatting, such as the date style on column D, is inherited from
# the row being pushed down, matching Excel's native Insert behaviour).
$ws.Rows(31).Insert()

# Populate the new row 31 with the new weekly record.
$ws.Cells.Item(31, 1).Value2  = 9
$ws.Cells.Item(31, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(31, 3).Value2  = "Metropolitana"
$ws.Cells.Item(31, 4).Value2  = 44742
$ws.Cells.Item(31, 5).Value2  = 13
$ws.Cells.Item(31, 6).Value2  = 100112022
$ws.Cells.Item(31, 7).Value2  = "Arveja Verde"
$ws.Cells.Item(31, 8).Value2  = "Perfection"
$ws.Cells.Item(31, 9).Value2  = "Primera"
$ws.Cells.Item(31, 10).Value2 = 43
$ws.Cells.Item(31, 11).Value2 = 38000
$ws.Cells.Item(31, 12).Value2 = 38000
$ws.Cells.Item(31, 13).Value2 = 38000
$ws.Cells.Item(31, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value2 = 1520
$ws.Cells.Item(31, 17).Value2 = 25
$ws.Cells.Item(31, 18).Value2 = "Hortaliza"
